$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-265)
# from serial date 45175 (2023-09-06) to serial date 45177 (2023-09-08)
for ($r = 2; $r -le 265; $r++) {
    $ws.Range("C$r").Value = 45177
}
